# Table.xlsx - mcc calc sheet
#
# Commit: "removed models from git and added them to gitignore.
#          Fixed spelling mistake in mcc calc sheet"
#
# The "removed models from git / .gitignore" part of the commit refers to
# other files in the repo (not this worksheet) and has no effect here.
#
# The spreadsheet-visible edit is the spelling fix: every occurrence of the
# mis-spelled row label "Sentsitivity" becomes "Sensitivity" (rows 2, 6, 10
# and 14 of the "Quality Measure" column, one per Outcome group). Excel's
# shared-string table naturally rebuilds as a result (the now-unused
# "Sentsitivity" entry drops out and a single new "Sensitivity" entry is
# appended at the end) - that is a mechanical side effect of the text edit,
# not a separate change.
#
# The author's selection also moved from H13 to D7 before saving, which we
# reproduce too.
#
# (Cosmetic deltas elsewhere in the saved XML - window size, the absPath
# hint, default row height / column "best fit" widths, x14ac:dyDescent -
# are artifacts of the file being opened/re-saved on a different PC/Excel
# build, per the changed "C:\Users\bobby\...\" -> "C:\Users\zephy\...\"
# path; they don't correspond to any deliberate user action and aren't
# reproducible user edits.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$correct = "Sensitivity"
$typo = "Sentsitivity"

# Note: read via .Value2 (.Value's getter misbehaves for string reads in
# this host); writes still go through .Value so cell type/shared-string
# handling stays exactly as a normal Excel edit would produce.
$used = $ws.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq $typo) {
        $cell.Value = $correct
    }
}

# Matches the final selection recorded in the saved file.
$ws.Range("D7").Select()
